$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cell F1 from "Type of hearing" to "Hearing type"
$ws.Range("F1").Value = "Hearing type"

# Update the active selection to match the edited cell
$ws.Range("F1").Select()
